# Automatic update of files.
# - Column C (Förändrad) bumps from 46073 to 46074 for every data row (2-14).
# - Rows 6,8,9,10,11,12,14 are re-shuffled: the A (Beteckning), B (Datum) and
#   G (Area (ha)) values move to different rows (two rotation cycles:
#   6->12->9->10->14->6 and 8->11->8), while everything else in those rows
#   (D, E, H..R) stays put because it is identical across the cycle.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C: Förändrad date bump (applies to every populated row 2-14) ---
foreach ($r in 2..14) {
    $ws.Cells.Item($r, 3).Value = 46074
}

# --- Rows 6, 8, 9, 10, 11, 12, 14: rotate Beteckning / Datum / Area values ---
# Capture the "before" values for A/B/G on the affected rows first, since the
# destination rows overlap with the source rows.
$rows = @(6, 8, 9, 10, 11, 12, 14)
$before = @{}
foreach ($r in $rows) {
    $before[$r] = @(
        $ws.Cells.Item($r, 1).Value2,   # A - Beteckning
        $ws.Cells.Item($r, 2).Value2,   # B - Datum
        $ws.Cells.Item($r, 7).Value2    # G - Area (ha)
    )
}

# Mapping of "data that used to live at row X now lives at row Y"
$moveTo = @{ 6 = 12; 8 = 11; 9 = 10; 10 = 14; 11 = 8; 12 = 9; 14 = 6 }

foreach ($srcRow in $rows) {
    $dstRow = $moveTo[$srcRow]
    $vals = $before[$srcRow]
    $ws.Cells.Item($dstRow, 1).Value = $vals[0]
    $ws.Cells.Item($dstRow, 2).Value = $vals[1]
    $ws.Cells.Item($dstRow, 7).Value = $vals[2]
}
